$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates per diff: N2 (text date) and numeric columns O2:AG2
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 996794422.3099999
$ws.Range("P2").Value = 388880204.53
$ws.Range("Q2").Value = 26819204.79
$ws.Range("R2").Value = -64.27829816800001
$ws.Range("S2").Value = 86257169.16
$ws.Range("T2").Value = 1.363159096
$ws.Range("U2").Value = 24590639.25
$ws.Range("V2").Value = 28.814932169
$ws.Range("W2").Value = 246481891.98
$ws.Range("X2").Value = 24174829.53
$ws.Range("Y2").Value = 53.2080219766
$ws.Range("Z2").Value = 18124808.24
$ws.Range("AA2").Value = 5.2656924779
$ws.Range("AB2").Value = 750312530.33
$ws.Range("AC2").Value = 12.6270181572
$ws.Range("AD2").Value = 13.8363894769
$ws.Range("AE2").Value = 17.6830895047
$ws.Range("AF2").Value = 212.5986447339
$ws.Range("AG2").Value = 24.727454976
